# TC07_Canine_Filter_SamplePatho-Osteosarcoma.xlsx
#
# The "startup" sheet carries, in column C (StatQuery), a Cypher query shared
# by the CasesTab / SamplesTab / FilesTab rows (C2:C4 all reference the same
# shared string). Two of its filter clauses need to swap which property is
# constrained to "Osteosarcoma":
#   - demo.breed filter goes from ['Osteosarcoma'] back to unrestricted []
#   - samp.specific_sample_pathology filter goes from unrestricted []
#     to ['Osteosarcoma']
#
# Also bumps the sheet's zoom level (cosmetic view state saved with the file).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldBreed = "AND (size(['Osteosarcoma']) = 0 OR demo.breed IN ['Osteosarcoma'])"
$newBreed = "AND (size([]) = 0 OR demo.breed IN [])"

$oldPath = "AND (size([]) = 0 OR samp.specific_sample_pathology IN [])"
$newPath = "AND (size(['Osteosarcoma']) = 0 OR samp.specific_sample_pathology IN ['Osteosarcoma'])"

foreach ($addr in @("C2", "C3", "C4")) {
    $cell = $ws.Range($addr)
    $text = $cell.Value2
    if ($text -ne $null) {
        $text = $text.Replace($oldBreed, $newBreed)
        $text = $text.Replace($oldPath, $newPath)
        $cell.Value2 = $text
    }
}

# Cosmetic: sheet zoom level changed from 55% to 25%.
$excel.ActiveWindow.Zoom = 25
